$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.16971755027771
$ws.Range("B1").Value = 2.43874979019165
$ws.Range("D1").Value = 2.366850852966309
$ws.Range("E1").Value = 1.233886003494263
